$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be introduced in this exact order to match
# the target shared-string table ordering (indices 31..42):
#  31 F6, 32 F7, 33 C6, 34 C7, 35 B6, 36 B7, 37 B8, 38 E8, 39 F8, 40 B9, 41 B13, 42 D10
$ws.Range("F6").Value = "Insert sound events for walking and"
$ws.Range("F7").Value = " background music on start"
$ws.Range("C6").Value = "Finish prototype segment: Level 1"
$ws.Range("C7").Value = "Start working in engine"
$ws.Range("B6").Value = "Start work in engine"
$ws.Range("B7").Value = "find more sound assets"
$ws.Range("B8").Value = "make UI assets (at least placeholder?)"
$ws.Range("E8").Value = "Cheat Codes: Reload Level shortcut, Load Next Level"
$ws.Range("F8").Value = "start Audio Guide"
$ws.Range("B9").Value = "Start Art Guide"
$ws.Range("B13").Value = "Art Guide"
$ws.Range("D10").Value = "Particle Effects"

# Row 9: A9 (Week 3 label) moves down to A10
$ws.Range("A9").ClearContents()
$ws.Range("A10").Value = "Week 3"

# Update the selected cell to D10
$ws.Range("D10").Select()
